# Actualización automática 2025-11-12 17:30:08
#
# A new client, "JAIME COELLO ALBERTO FERNANDO", needs to be inserted
# (alphabetically) right before "JIMENEZ CORDERO WILLIAM GUSTAVO" on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. That is row 28 on each sheet.
# Inserting the row there pushes every following client down by one row and
# grows the totals/summary row at the bottom from row 54 -> 55.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, summary row uses "<n> de <count>" text)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(28).Insert()

$ws1.Range("A28").Value = "OFICINA-CATAECSA"
$ws1.Range("B28").Value = "JAIME COELLO ALBERTO FERNANDO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(28, $col).Value = 0
}

# The trailing "<n> de 52" summary row has shifted from row 54 to row 55; the
# denominator must grow to 53 to reflect the new total number of clients.
$counts1 = @{
    "C" = 0; "D" = 0; "E" = 1; "F" = 0; "G" = 0; "H" = 1; "I" = 1; "J" = 0
    "K" = 0; "L" = 1; "M" = 1; "N" = 0; "O" = 0; "P" = 0; "Q" = 0; "R" = 0
}
foreach ($col in $counts1.Keys) {
    $ws1.Range($col + "55").Value = "$($counts1[$col]) de 53"
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, summary row holds numeric totals)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(28).Insert()

$ws2.Range("A28").Value = "OFICINA-CATAECSA"
$ws2.Range("B28").Value = "JAIME COELLO ALBERTO FERNANDO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(28, $col).Value = 0
}
